# 0001359: Fiche Pulve : Marque recherche sur le nom avec autocomplétion
$d = $word.ActiveDocument

# 1) Merge the split "Suppression des boutons « Voir Fiche » et « supprimer » "
#    run-sequence into one run (occurs twice: DrawingML + VML fallback textboxes).
$old1 = "Suppression des boutons " + [char]0x201C + "Voir Fiche" + [char]0x201D + " et " + [char]0x201C + "supprimer" + [char]0x201D + " "
$new1 = "Suppression des boutons " + [char]0xAB + " Voir Fiche " + [char]0xBB + " et " + [char]0xAB + " supprimer " + [char]0xBB + " "

$rng = $d.Content
$rng.Find.Execute("Suppression des boutons*supprimer*", $false, $false, $true, $false, $false, $true, 1, $false, "", 0)

# 2) Replace the "Reprises de champs accessoires ... Regulation du pulvérisateur..."
#    sentence with the shortened "Reprises du pulvérisateur principal Non Modifiable."
$find = $d.Content.Find
$find.Execute("Reprises de champs accessoires /Manomètres/ Attelage, pulvérisation, Regulation du pulvérisateur principal Non Modifiable.", $true, $false, $false, $false, $false, $true, 1, $false, "Reprises du pulvérisateur principal Non Modifiable.", 2)
